# Update "paises" (countries) COVID dashboard snapshot.
# - Refreshes the "Datos actualizados..." timestamp cell.
# - Swaps the country names for four pairs of adjacent-rank rows whose case
#   counts crossed each other (Argentina/Arabia Saudita, Japon/Singapur,
#   Surinam/Mayotte, Montserrat/Islas Malvinas).
# - Refreshes the numeric columns (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Muertes hoy, Muertes) for the rows whose figures
#   changed in this data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 01:54"

# --- Country-name swaps (rank changed) ---------------------------------
$ws.Range("A16").Value = "Argentina"
$ws.Range("A17").Value = "Arabia Saudita"

$ws.Range("A48").Value = "Japon"
$ws.Range("A49").Value = "Singapur"

$ws.Range("A121").Value = "Surinam"
$ws.Range("A122").Value = "Mayotte"

$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Numeric refresh -----------------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5653293
$ws.Range("C4").Value = 41318
$ws.Range("D4").Value = 3002504
$ws.Range("E4").Value = 2475848
$ws.Range("G4").Value = 1225
$ws.Range("H4").Value = 174941

# Row 5 - Brasil
$ws.Range("B5").Value = 3411872
$ws.Range("C5").Value = 48637
$ws.Range("E5").Value = 747674
$ws.Range("G5").Value = 1365
$ws.Range("H5").Value = 110019

# Row 6 - India
$ws.Range("D6").Value = 2036703
$ws.Range("E6").Value = 676900

# Row 9 - Peru
$ws.Range("B9").Value = 549321
$ws.Range("C9").Value = 7828
$ws.Range("D9").Value = 374019
$ws.Range("E9").Value = 148644
$ws.Range("G9").Value = 177
$ws.Range("H9").Value = 26658

# Row 16 - now Argentina
$ws.Range("B16").Value = 305966
$ws.Range("C16").Value = 6840
$ws.Range("D16").Value = 223531
$ws.Range("E16").Value = 76387
$ws.Range("G16").Value = 234
$ws.Range("H16").Value = 6048

# Row 17 - now Arabia Saudita
$ws.Range("B17").Value = 301323
$ws.Range("C17").Value = 1409
$ws.Range("D17").Value = 272911
$ws.Range("E17").Value = 24942
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 3470

# Row 23 - Francia
$ws.Range("E23").Value = 106751
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 30451

# Row 27 - Canada
$ws.Range("B27").Value = 123154
$ws.Range("C27").Value = 282
$ws.Range("D27").Value = 109357
$ws.Range("E27").Value = 4752

# Row 48 - now Japon
$ws.Range("B48").Value = 56685
$ws.Range("C48").Value = 1018
$ws.Range("D48").Value = 42284
$ws.Range("E48").Value = 13286
$ws.Range("G48").Value = 16
$ws.Range("H48").Value = 1115

# Row 49 - now Singapur
$ws.Range("B49").Value = 55938
$ws.Range("C49").Value = 100
$ws.Range("D49").Value = 52533
$ws.Range("E49").Value = 3378
$ws.Range("H49").Value = 27

# Row 53 - Barein
$ws.Range("B53").Value = 47581
$ws.Range("C53").Value = 396
$ws.Range("D53").Value = 43921
$ws.Range("E53").Value = 3485

# Row 74 - Chequia
$ws.Range("B74").Value = 20483
$ws.Range("C74").Value = 281
$ws.Range("E74").Value = 4936

# Row 113 - Nicaragua
$ws.Range("B113").Value = 4311
$ws.Range("C113").Value = 196
$ws.Range("E113").Value = 1265
$ws.Range("G113").Value = 5
$ws.Range("H113").Value = 133

# Row 114 - Montenegro
$ws.Range("B114").Value = 4132
$ws.Range("C114").Value = 47
$ws.Range("D114").Value = 3035
$ws.Range("E114").Value = 1017

# Row 121 - now Surinam
$ws.Range("B121").Value = 3216
$ws.Range("C121").Value = 139
$ws.Range("D121").Value = 2196
$ws.Range("E121").Value = 966
$ws.Range("G121").Value = 6
$ws.Range("H121").Value = 54

# Row 122 - now Mayotte
$ws.Range("B122").Value = 3160
$ws.Range("D122").Value = 2964
$ws.Range("E122").Value = 157
$ws.Range("H122").Value = 39

# Row 143 - Uruguay
$ws.Range("B143").Value = 1485
$ws.Range("C143").Value = 28
$ws.Range("D143").Value = 1219
$ws.Range("E143").Value = 226

# Row 146 - Republica de Chipre
$ws.Range("B146").Value = 1359
$ws.Range("C146").Value = 8
$ws.Range("D146").Value = 878

# Row 167 - Guadalupe
$ws.Range("E167").Value = 206
$ws.Range("H167").Value = 15

# Row 213 - now Montserrat
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 - now Islas Malvinas
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
